$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a text value that looks like a percentage (e.g. "68%")
# without letting Excel auto-convert it to a numeric percent value,
# and without leaving the cell on a different style record than the
# rest of the column (style donor = H2, never itself modified).
function Set-PercentText($ref, $value) {
    $cell = $ws.Range($ref)
    $donor = $ws.Range("H2")
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $donor.Copy() | Out-Null
    $cell.PasteSpecial(-4122) | Out-Null
}

$ws.Range("E2").Value = '2026-02-06 15:47:44'
$ws.Range("I2").Value = '0.3 mm'
$ws.Range("K2").Value = '8.5 MJ/m2'
$ws.Range("O2").Value = '-0.2 °C'
$ws.Range("E3").Value = '2026-02-06 15:47:46'
$ws.Range("K3").Value = '12.0 MJ/m2'
$ws.Range("E4").Value = '2026-02-06 15:47:48'
$ws.Range("J4").Value = '996.7 hPa'
$ws.Range("K4").Value = '11.1 MJ/m2'
$ws.Range("O4").Value = '13.5 °C'
$ws.Range("E5").Value = '2026-02-06 15:47:51'
Set-PercentText "H5" '68%'
$ws.Range("J5").Value = '997.0 hPa'
$ws.Range("K5").Value = '10.2 MJ/m2'
$ws.Range("O5").Value = '10.8 °C'
$ws.Range("E6").Value = '2026-02-06 15:47:53'
$ws.Range("J6").Value = '998.2 hPa'
$ws.Range("K6").Value = '9.4 MJ/m2'
$ws.Range("O6").Value = '15.4 °C'
$ws.Range("E7").Value = '2026-02-06 15:47:56'
$ws.Range("J7").Value = '997.8 hPa'
$ws.Range("K7").Value = '11.5 MJ/m2'
$ws.Range("E8").Value = '2026-02-06 15:47:58'
Set-PercentText "H8" '76%'
$ws.Range("K8").Value = '11.6 MJ/m2'
$ws.Range("O8").Value = '10.0 °C'
$ws.Range("E9").Value = '2026-02-06 15:48:00'
Set-PercentText "H9" '85%'
$ws.Range("O9").Value = '4.6 °C'
$ws.Range("E10").Value = '2026-02-06 15:48:03'
$ws.Range("H10").ClearContents()
$ws.Range("I10").ClearContents()
$ws.Range("M10").ClearContents()
$ws.Range("N10").ClearContents()
$ws.Range("O10").ClearContents()
$ws.Range("E11").Value = '2026-02-06 15:48:04'
Set-PercentText "H11" '77%'
$ws.Range("K11").Value = '8.6 MJ/m2'
$ws.Range("O11").Value = '5.2 °C'
$ws.Range("E12").Value = '2026-02-06 15:48:07'
$ws.Range("K12").Value = '11.8 MJ/m2'
$ws.Range("O12").Value = '14.4 °C'
$ws.Range("E13").Value = '2026-02-06 15:48:09'
Set-PercentText "H13" '76%'
$ws.Range("O13").Value = '10.0 °C'
$ws.Range("E14").Value = '2026-02-06 15:48:11'
Set-PercentText "H14" '72%'
$ws.Range("K14").Value = '7.2 MJ/m2'
$ws.Range("E15").Value = '2026-02-06 15:48:14'
Set-PercentText "H15" '72%'
$ws.Range("J15").Value = '997.1 hPa'
$ws.Range("K15").Value = '11.4 MJ/m2'
$ws.Range("O15").Value = '10.2 °C'
$ws.Range("E16").Value = '2026-02-06 15:48:16'
Set-PercentText "H16" '86%'
$ws.Range("K16").Value = '9.3 MJ/m2'
$ws.Range("L16").Value = '25.6 km/h - 218º 15:26 TU'
$ws.Range("O16").Value = '5.8 °C'
$ws.Range("E17").Value = '2026-02-06 15:48:19'
Set-PercentText "H17" '85%'
$ws.Range("K17").Value = '10.3 MJ/m2'
$ws.Range("O17").Value = '5.7 °C'
$ws.Range("E18").Value = '2026-02-06 15:48:21'
$ws.Range("K18").Value = '5.7 MJ/m2'
$ws.Range("E19").Value = '2026-02-06 15:48:24'
Set-PercentText "H19" '77%'
$ws.Range("J19").Value = '999.5 hPa'
$ws.Range("K19").Value = '11.4 MJ/m2'
$ws.Range("O19").Value = '9.7 °C'
$ws.Range("E20").Value = '2026-02-06 15:48:26'
$ws.Range("K20").Value = '11.7 MJ/m2'
$ws.Range("O20").Value = '-1.8 °C'
$ws.Range("E21").Value = '2026-02-06 15:48:29'
Set-PercentText "H21" '74%'
$ws.Range("J21").Value = '997.4 hPa'
$ws.Range("K21").Value = '10.0 MJ/m2'
$ws.Range("O21").Value = '8.2 °C'
$ws.Range("E22").Value = '2026-02-06 15:48:31'
$ws.Range("K22").Value = '11.3 MJ/m2'
$ws.Range("O22").Value = '10.6 °C'
$ws.Range("E23").Value = '2026-02-06 15:48:33'
Set-PercentText "H23" '81%'
$ws.Range("J23").Value = '997.1 hPa'
$ws.Range("K23").Value = '8.7 MJ/m2'
$ws.Range("O23").Value = '9.9 °C'
$ws.Range("E24").Value = '2026-02-06 15:48:36'
$ws.Range("J24").Value = '996.5 hPa'
$ws.Range("K24").Value = '11.5 MJ/m2'
$ws.Range("E25").Value = '2026-02-06 15:48:38'
Set-PercentText "H25" '80%'
$ws.Range("K25").Value = '9.3 MJ/m2'
$ws.Range("O25").Value = '4.1 °C'
$ws.Range("E26").Value = '2026-02-06 15:48:41'
$ws.Range("K26").Value = '7.8 MJ/m2'
$ws.Range("E27").Value = '2026-02-06 15:48:43'
Set-PercentText "H27" '82%'
$ws.Range("J27").Value = '997.1 hPa'
$ws.Range("K27").Value = '10.0 MJ/m2'
$ws.Range("O27").Value = '10.6 °C'
$ws.Range("E28").Value = '2026-02-06 15:48:46'
Set-PercentText "H28" '84%'
$ws.Range("O28").Value = '4.5 °C'
$ws.Range("E29").Value = '2026-02-06 15:48:48'
$ws.Range("K29").Value = '11.6 MJ/m2'
$ws.Range("O29").Value = '12.6 °C'
$ws.Range("E30").Value = '2026-02-06 15:48:50'
$ws.Range("K30").Value = '8.8 MJ/m2'
$ws.Range("E31").Value = '2026-02-06 15:48:53'
Set-PercentText "H31" '85%'
$ws.Range("J31").Value = '998.8 hPa'
$ws.Range("O31").Value = '7.1 °C'
$ws.Range("E32").Value = '2026-02-06 15:48:55'
$ws.Range("J32").Value = '998.5 hPa'
$ws.Range("K32").Value = '11.8 MJ/m2'
$ws.Range("E33").Value = '2026-02-06 15:48:57'
$ws.Range("O33").Value = '10.0 °C'
$ws.Range("E34").Value = '2026-02-06 15:49:00'
Set-PercentText "H34" '76%'
$ws.Range("K34").Value = '11.5 MJ/m2'
$ws.Range("O34").Value = '8.5 °C'
$ws.Range("E35").Value = '2026-02-06 15:49:02'
$ws.Range("K35").Value = '9.2 MJ/m2'
$ws.Range("E36").Value = '2026-02-06 15:49:04'
Set-PercentText "H36" '59%'
$ws.Range("J36").Value = '999.5 hPa'
$ws.Range("K36").Value = '11.1 MJ/m2'
$ws.Range("O36").Value = '13.3 °C'
